# Adds attendance rows 7-21 (days 2018-01-02 .. 01-04), mirroring rows 2-6's
# formatting: col A = date style, B/C/D/E = time style, F = plain weekday code.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A2:C2").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Range("A7").Value = 43102
$ws.Range("B7").Value = 0.3263888888888889
$ws.Range("C7").Value = 0.70833333333333337
$ws.Range("F7").Value = 1

# Row 8
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 43102
$ws.Range("F8").Value = 2

# Row 9
$ws.Range("A2:C2").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Range("A9").Value = 43102
$ws.Range("B9").Value = 0.3263888888888889
$ws.Range("C9").Value = 0.70833333333333337
$ws.Range("F9").Value = 4

# Row 10
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 43102
$ws.Range("F10").Value = 5

# Row 11
$ws.Range("A2:C2").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A11").Value = 43102
$ws.Range("B11").Value = 0.3263888888888889
$ws.Range("C11").Value = 0.70833333333333337
$ws.Range("F11").Value = 6

# Row 12
$ws.Range("A2:C2").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A12").Value = 43103
$ws.Range("B12").Value = 0.3263888888888889
$ws.Range("C12").Value = 0.70833333333333337
$ws.Range("F12").Value = 1

# Row 13
$ws.Range("A2:C2").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A13").Value = 43103
$ws.Range("B13").Value = 0.32291666666666669
$ws.Range("C13").Value = 0.70833333333333304
$ws.Range("F13").Value = 2

# Row 14
$ws.Range("A3:E3").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A14").Value = 43103
$ws.Range("B14").Value = 0.32291666666666669
$ws.Range("C14").Value = 0.41666666666666669
$ws.Range("D14").Value = 0.54166666666666663
$ws.Range("E14").Value = 0.70833333333333337
$ws.Range("F14").Value = 4

# Row 15
$ws.Range("A2:C2").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = 43103
$ws.Range("B15").Value = 0.3263888888888889
$ws.Range("C15").Value = 0.70833333333333304
$ws.Range("F15").Value = 5

# Row 16
$ws.Range("A2:C2").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = 43103
$ws.Range("B16").Value = 0.32291666666666669
$ws.Range("C16").Value = 0.70833333333333304
$ws.Range("F16").Value = 6

# Row 17
$ws.Range("A3:E3").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("A17").Value = 43104
$ws.Range("B17").Value = 0.3263888888888889
$ws.Range("C17").Value = 0.41666666666666669
$ws.Range("D17").Value = 0.54166666666666663
$ws.Range("E17").Value = 0.70833333333333337
$ws.Range("F17").Value = 1

# Row 18
$ws.Range("A2:C2").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A18").Value = 43104
$ws.Range("B18").Value = 0.33680555555555558
$ws.Range("C18").Value = 0.70833333333333304
$ws.Range("F18").Value = 2

# Row 19
$ws.Range("A2:C2").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A19").Value = 43104
$ws.Range("B19").Value = 0.3263888888888889
$ws.Range("C19").Value = 0.70833333333333304
$ws.Range("F19").Value = 4

# Row 20
$ws.Range("A2:C2").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = 43104
$ws.Range("B20").Value = 0.3263888888888889
$ws.Range("C20").Value = 0.70833333333333304
$ws.Range("F20").Value = 5

# Row 21
$ws.Range("A2:C2").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = 43104
$ws.Range("B21").Value = 0.32291666666666669
$ws.Range("C21").Value = 0.70833333333333304
$ws.Range("F21").Value = 6

$excel.CutCopyMode = $false

# Final selection, matching the saved workbook state
[void]$ws.Range("H9").Select()
